$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 59
$ws.Range("F8").Value = 30
$ws.Range("F9").Value = 280
$ws.Range("F10").Value = 405
$ws.Range("F12").Value = 802
$ws.Range("F13").Value = 790
$ws.Range("F14").Value = 10
$ws.Range("F15").Value = 12
$ws.Range("F16").Value = 1552
$ws.Range("F17").Value = 1552
$ws.Range("F18").Value = 1192
$ws.Range("F21").Value = 169
$ws.Range("F22").Value = 371
$ws.Range("F25").Value = 118
$ws.Range("F26").Value = 6774
$ws.Range("F27").Value = 5202
$ws.Range("F28").Value = 12
$ws.Range("F29").Value = 156
$ws.Range("F32").Value = 219
$ws.Range("F37").Value = 1327
$ws.Range("F39").Value = 264
$ws.Range("F40").Value = 636
$ws.Range("F43").Value = 275
$ws.Range("F45").Value = 159
$ws.Range("F47").Value = 93
$ws.Range("F48").Value = 107

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 43
$ws.Range("F10").Value = 17
$ws.Range("F15").Value = 54
$ws.Range("F18").Value = 255

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2484
$ws.Range("F4").Value = 222
$ws.Range("F5").Value = 86

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 222
$ws.Range("F7").Value = 86
$ws.Range("F10").Value = 59
$ws.Range("F11").Value = 30
$ws.Range("F12").Value = 280
$ws.Range("F14").Value = 405
$ws.Range("F16").Value = 802
$ws.Range("F17").Value = 790
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 12
$ws.Range("F20").Value = 1552
$ws.Range("F21").Value = 1552
$ws.Range("F22").Value = 1192
$ws.Range("F24").Value = 371
$ws.Range("F26").Value = 118
$ws.Range("F27").Value = 43
$ws.Range("F29").Value = 6774
$ws.Range("F30").Value = 5202
$ws.Range("F31").Value = 219
$ws.Range("F33").Value = 1327
$ws.Range("F36").Value = 264
$ws.Range("F38").Value = 636
$ws.Range("F41").Value = 54
$ws.Range("F43").Value = 275
$ws.Range("F44").Value = 159
$ws.Range("F46").Value = 93
$ws.Range("F47").Value = 107
$ws.Range("F49").Value = 255
